$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for refreshed values ---
$ws.Range("D2").Value = "43.647.05"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "2.370.66"
$ws.Range("E3").Value = "  +5.91%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'234.95"
$ws.Range("E5").Value = "  +2.12%  "

$ws.Range("D6").Value = "'0.655"
$ws.Range("E6").Value = "  +2.48%  "

$ws.Range("D7").Value = "'73.69"
$ws.Range("E7").Value = "  +15.52%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +17.94%  "

$ws.Range("D10").Value = "'0.0984"
$ws.Range("E10").Value = "  +3.84%  "

$ws.Range("D11").Value = "'27.75"
$ws.Range("E11").Value = "  +4.19%  "

$ws.Range("D12").Value = "2.717.47"
$ws.Range("E12").Value = "  +5.62%  "

$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").Value = "'16.49"
$ws.Range("E14").Value = "  +8.90%  "

$ws.Range("D15").Value = "'6.35"
$ws.Range("E15").Value = "  +6.03%  "

$ws.Range("D16").Value = "'0.876"
$ws.Range("E16").Value = "  +6.71%  "

$ws.Range("D17").Value = "2.370.24"
$ws.Range("E17").Value = "  +5.80%  "

$ws.Range("D18").Value = "43.512.71"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("E19").Value = "  +5.28%  "

$ws.Range("D20").Value = "'75.53"
$ws.Range("E20").Value = "  +3.69%  "

$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  +6.31%  "

$ws.Range("D22").Value = "'252.04"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +3.66%  "

$ws.Range("D26").Value = "'10.18"
$ws.Range("E26").Value = "  +5.26%  "

$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("D28").Value = "'22.65"
$ws.Range("E28").Value = "  +5.12%  "

$ws.Range("D29").Value = "'172.76"
$ws.Range("E29").Value = "  -0.81%  "

$ws.Range("E30").Value = "  +9.74%  "

$ws.Range("E31").Value = "  +4.05%  "

$ws.Range("D32").Value = "'0.128"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("E33").Value = "  +3.83%  "

$ws.Range("D34").Value = "'0.0703"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("D35").Value = "'5.12"
$ws.Range("E35").Value = "  +4.90%  "

$ws.Range("D36").Value = "'3.79"
$ws.Range("E36").Value = "  +6.31%  "

$ws.Range("D37").Value = "'6.67"
$ws.Range("E37").Value = "  +6.07%  "

$ws.Range("D38").Value = "'2.46"
$ws.Range("E38").Value = "  +9.30%  "

$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("D40").Value = "'19.46"
$ws.Range("E40").Value = "  +15.18%  "

$ws.Range("E41").Value = "  +4.48%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").Value = "'100.60"
$ws.Range("E43").Value = "  +4.69%  "

$ws.Range("D46").Value = "'1.22"
$ws.Range("E46").Value = "  +3.44%  "

$ws.Range("D47").Value = "'0.0961"
$ws.Range("E47").Value = "  +2.43%  "

$ws.Range("D48").Value = "1.452.10"
$ws.Range("E48").Value = "  +1.84%  "

# --- Row 44/45 swap (FTXToken <-> ARBITRUM) with refreshed values ---
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  +11.08%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.52"
$ws.Range("E45").Value = "  +1.49%  "

# --- Row 49/50 swap (RocketPoolETH <-> Algorand) with refreshed values ---
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.175"
$ws.Range("E49").Value = "  +9.22%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.594.75"
$ws.Range("E50").Value = "  +6.00%  "

# --- Row 51: HuobiToken -> TerraClassic ---
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "'0.000203"
$ws.Range("E51").Value = "  -1.60%  "
